$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.04
$ws.Range("G2").Value = 1000
$ws.Range("H2").Value = 1.04
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 1.01
$ws.Range("K2").Value = 950
$ws.Range("P2").Value = 1.24
$ws.Range("Q2").Value = 1.01
# Row 3
$ws.Range("F3").Value = 1.04
$ws.Range("G3").Value = 970
$ws.Range("H3").Value = 1.04
$ws.Range("I3").Value = 970
$ws.Range("J3").Value = 1.03
$ws.Range("K3").Value = 950
$ws.Range("P3").Value = 1.25
$ws.Range("Q3").Value = 1.25
$ws.Range("R3").Value = 1.18
$ws.Range("S3").Value = 1.25
# Row 4
$ws.Range("F4").Value = 1.04
$ws.Range("H4").Value = 1.04
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 1.01
$ws.Range("N4").Value = 1.01
$ws.Range("P4").Value = 1.25
$ws.Range("R4").Value = 1.18
$ws.Range("S4").Value = 1.09
$ws.Range("V4").Value = 1.01
# Row 5
$ws.Range("J5").Value = 3.95
$ws.Range("O5").Value = 1.18
$ws.Range("T5").Value = 1.51
$ws.Range("X5").Value = 1000
$ws.Range("Y5").Value = 18
$ws.Range("Z5").Value = 980
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 1000
$ws.Range("AD5").Value = 12.5
$ws.Range("AE5").Value = 1000
$ws.Range("AF5").Value = 1000
$ws.Range("AG5").Value = 18
$ws.Range("AH5").Value = 980
$ws.Range("AI5").Value = 1000
$ws.Range("AJ5").Value = 1000
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 1000
$ws.Range("AO5").Value = 11.5
# Row 6
$ws.Range("F6").Value = 3.3
$ws.Range("G6").Value = 3.45
$ws.Range("H6").Value = 2.16
$ws.Range("I6").Value = 2.18
$ws.Range("K6").Value = 4.2
$ws.Range("N6").Value = 5.6
$ws.Range("O6").Value = 1.16
$ws.Range("P6").Value = 2.74
$ws.Range("R6").Value = 1.69
$ws.Range("T6").Value = 1.49
$ws.Range("X6").Value = 980
$ws.Range("Y6").Value = 20
$ws.Range("Z6").Value = 980
$ws.Range("AA6").Value = 980
$ws.Range("AB6").Value = 980
$ws.Range("AD6").Value = 13
$ws.Range("AE6").Value = 980
$ws.Range("AF6").Value = 980
$ws.Range("AG6").Value = 18
$ws.Range("AI6").Value = 980
$ws.Range("AK6").Value = 1000
$ws.Range("AL6").Value = 1000
$ws.Range("AN6").Value = 980
$ws.Range("AO6").Value = 11
# Row 8
$ws.Range("F8").Value = 1.04
$ws.Range("G8").Value = 980
$ws.Range("H8").Value = 1.29
$ws.Range("I8").Value = 980
$ws.Range("J8").Value = 1.01
$ws.Range("K8").Value = 980
$ws.Range("P8").Value = 1.16
$ws.Range("Q8").Value = 1.01
# Row 9
$ws.Range("F9").Value = 1.04
$ws.Range("G9").Value = 980
$ws.Range("H9").Value = 1.04
$ws.Range("J9").Value = 1.02
$ws.Range("K9").Value = 980
$ws.Range("P9").Value = 1.15
$ws.Range("Q9").Value = 1.01
# Row 10
$ws.Range("F10").Value = 1.04
$ws.Range("G10").Value = 980
$ws.Range("H10").Value = 1.04
$ws.Range("I10").Value = 980
$ws.Range("J10").Value = 1.02
$ws.Range("K10").Value = 980
$ws.Range("P10").Value = 1.15
$ws.Range("Q10").Value = 1.01
# Row 11
$ws.Range("F11").Value = 1.04
$ws.Range("G11").Value = 980
$ws.Range("H11").Value = 1.04
$ws.Range("I11").Value = 980
$ws.Range("J11").Value = 1.02
$ws.Range("K11").Value = 980
$ws.Range("P11").Value = 1.16
$ws.Range("Q11").Value = 1.01
# Row 12
$ws.Range("G12").Value = 1.87
$ws.Range("H12").Value = 4.7
$ws.Range("O12").Value = 1.32
$ws.Range("P12").Value = 1.91
$ws.Range("S12").Value = 3.6
$ws.Range("U12").Value = 2
$ws.Range("X12").Value = 15
$ws.Range("Y12").Value = 1000
$ws.Range("Z12").Value = 1000
$ws.Range("AA12").Value = 1000
$ws.Range("AB12").Value = 9.800000000000001
$ws.Range("AC12").Value = 9.199999999999999
$ws.Range("AD12").Value = 1000
$ws.Range("AE12").Value = 1000
$ws.Range("AF12").Value = 12
$ws.Range("AG12").Value = 12
$ws.Range("AH12").Value = 1000
$ws.Range("AI12").Value = 1000
$ws.Range("AJ12").Value = 1000
$ws.Range("AK12").Value = 1000
$ws.Range("AL12").Value = 1000
$ws.Range("AM12").Value = 1000
$ws.Range("AN12").Value = 1000
$ws.Range("AO12").Value = 1000
# Row 13
$ws.Range("F13").Value = 8.800000000000001
$ws.Range("G13").Value = 9.6
$ws.Range("H13").Value = 1.34
$ws.Range("I13").Value = 1.36
$ws.Range("J13").Value = 6.6
$ws.Range("K13").Value = 7
$ws.Range("O13").Value = 1.11
$ws.Range("P13").Value = 3.25
$ws.Range("Q13").Value = 1.37
$ws.Range("R13").Value = 1.97
$ws.Range("S13").Value = 1.94
$ws.Range("T13").Value = 1.62
$ws.Range("U13").Value = 2.32
$ws.Range("X13").Value = 1000
$ws.Range("Y13").Value = 1000
$ws.Range("Z13").Value = 13
$ws.Range("AA13").Value = 1000
$ws.Range("AB13").Value = 1000
$ws.Range("AC13").Value = 1000
$ws.Range("AD13").Value = 12
$ws.Range("AE13").Value = 14.5
$ws.Range("AF13").Value = 1000
$ws.Range("AG13").Value = 1000
$ws.Range("AH13").Value = 1000
$ws.Range("AI13").Value = 1000
$ws.Range("AL13").Value = 1000
$ws.Range("AN13").Value = 1000
$ws.Range("AO13").Value = 3.85
# Row 14
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 4.1
$ws.Range("R14").Value = 1.52
$ws.Range("S14").Value = 2.7
$ws.Range("T14").Value = 1.61
$ws.Range("U14").Value = 2.42
$ws.Range("X14").Value = 29
$ws.Range("Y14").Value = 13.5
$ws.Range("Z14").Value = 17.5
$ws.Range("AA14").Value = 980
$ws.Range("AB14").Value = 980
$ws.Range("AC14").Value = 9.4
$ws.Range("AE14").Value = 980
$ws.Range("AF14").Value = 980
$ws.Range("AG14").Value = 980
$ws.Range("AH14").Value = 980
$ws.Range("AI14").Value = 980
$ws.Range("AJ14").Value = 1000
$ws.Range("AK14").Value = 980
$ws.Range("AM14").Value = 1000
$ws.Range("AN14").Value = 980
$ws.Range("AO14").Value = 12
# Row 15
$ws.Range("F15").Value = 1.9
$ws.Range("G15").Value = 1.94
$ws.Range("I15").Value = 4.8
$ws.Range("J15").Value = 3.7
$ws.Range("N15").Value = 3.6
$ws.Range("O15").Value = 1.35
$ws.Range("P15").Value = 1.88
$ws.Range("Q15").Value = 2.04
$ws.Range("R15").Value = 1.33
$ws.Range("S15").Value = 3.7
$ws.Range("T15").Value = 1.89
$ws.Range("U15").Value = 1.94
$ws.Range("X15").Value = 16
$ws.Range("Y15").Value = 980
$ws.Range("Z15").Value = 980
$ws.Range("AA15").Value = 1000
$ws.Range("AB15").Value = 9.6
$ws.Range("AC15").Value = 8.199999999999999
$ws.Range("AD15").Value = 980
$ws.Range("AE15").Value = 1000
$ws.Range("AF15").Value = 11.5
$ws.Range("AG15").Value = 12
$ws.Range("AH15").Value = 980
$ws.Range("AI15").Value = 1000
$ws.Range("AJ15").Value = 980
$ws.Range("AK15").Value = 980
$ws.Range("AL15").Value = 980
$ws.Range("AM15").Value = 1000
$ws.Range("AN15").Value = 16.5
# Row 16
$ws.Range("J16").Value = 5.5
$ws.Range("K16").Value = 5.7
$ws.Range("N16").Value = 7.4
$ws.Range("Q16").Value = 1.41
$ws.Range("T16").Value = 1.55
$ws.Range("U16").Value = 2.54
$ws.Range("X16").Value = 1000
$ws.Range("Y16").Value = 1000
$ws.Range("Z16").Value = 1000
$ws.Range("AB16").Value = 1000
$ws.Range("AC16").Value = 16
$ws.Range("AD16").Value = 1000
$ws.Range("AE16").Value = 1000
$ws.Range("AF16").Value = 14.5
$ws.Range("AH16").Value = 1000
$ws.Range("AI16").Value = 1000
$ws.Range("AJ16").Value = 19
$ws.Range("AK16").Value = 16
$ws.Range("AL16").Value = 1000
$ws.Range("AM16").Value = 1000
$ws.Range("AN16").Value = 4.7
$ws.Range("AO16").Value = 1000
# Row 17
$ws.Range("H17").Value = 21
$ws.Range("I17").Value = 25
$ws.Range("J17").Value = 8.6
$ws.Range("K17").Value = 9
$ws.Range("N17").Value = 8.199999999999999
$ws.Range("O17").Value = 1.1
$ws.Range("P17").Value = 3.25
$ws.Range("R17").Value = 1.97
$ws.Range("S17").Value = 1.94
$ws.Range("T17").Value = 2.08
$ws.Range("U17").Value = 1.77
$ws.Range("X17").Value = 1000
$ws.Range("Y17").Value = 1000
$ws.Range("AB17").Value = 980
$ws.Range("AC17").Value = 1000
$ws.Range("AF17").Value = 10
$ws.Range("AG17").Value = 980
$ws.Range("AH17").Value = 1000
$ws.Range("AJ17").Value = 9.6
$ws.Range("AK17").Value = 16
$ws.Range("AL17").Value = 980
# Row 18
$ws.Range("G18").Value = 1.14
$ws.Range("H18").Value = 24
$ws.Range("I18").Value = 26
$ws.Range("N18").Value = 9.6
$ws.Range("P18").Value = 4.2
$ws.Range("Q18").Value = 1.24
$ws.Range("T18").Value = 2.02
$ws.Range("AB18").Value = 1000
$ws.Range("AC18").Value = 1000
$ws.Range("AF18").Value = 12.5
$ws.Range("AG18").Value = 1000
$ws.Range("AJ18").Value = 10.5
$ws.Range("AK18").Value = 15.5

Write-Host "Applied 259 cell updates"